$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENTG")

# Row 4 - Inventory
$ws.Range("B4").Value = 359000000.0
$ws.Range("C4").Value = 324000000.0
$ws.Range("D4").Value = 330000000.0
$ws.Range("E4").Value = 333000000.0
$ws.Range("F4").Value = 301000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 93000000.0
$ws.Range("C15").Value = 82000000.0
$ws.Range("D15").Value = 70000000.0
$ws.Range("E15").Value = 80000000.0
$ws.Range("F15").Value = 82000000.0

# Row 24 - Long Term Tax Liability (Deferred)
$ws.Range("B24").Value = 58000000.0
$ws.Range("C24").Value = 59000000.0
$ws.Range("D24").Value = 76000000.0
$ws.Range("E24").Value = 60000000.0
$ws.Range("F24").Value = 60000000.0

# Row 37 - Net Debt
$ws.Range("G37").Value = 584573000.0

# Row 38 - Total Debt
$ws.Range("G38").Value = 936484000.0
